$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new dSF (column F) value, per repulled/pushed data.
$updates = @{
    3  = 5
    4  = -3
    5  = 2
    6  = 1
    7  = -1
    8  = 6
    9  = 4
    10 = 5
    11 = 7
    12 = -5
    13 = 1
    14 = -1
    15 = -1
    16 = 1
    17 = -1
    18 = 0
    19 = 1
    20 = 1
    21 = -2
    22 = 1
    23 = -1
    24 = -3
    25 = -2
    26 = -1
    27 = 5
    29 = -3
    30 = -2
    31 = 1
    32 = -1
    33 = 4
    34 = 2
    35 = 3
    36 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
